$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("I11").Value = 890.8
$ws1.Range("L22").Value = 1710.72
$ws1.Range("M45").Value = 4224.69
$ws1.Range("I56").Value = "6 de 54"
$ws1.Range("L56").Value = "10 de 54"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F11").Value = 890.8
$ws2.Range("F22").Value = 1710.72
$ws2.Range("F45").Value = 5256.7
$ws2.Range("F56").Value = 78869.53999999999

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D8").Value = 1859.49
$ws3.Range("E8").Value = -859.49
$ws3.Range("F8").Value = 1.85949

$ws3.Range("D15").Value = 16338.05
$ws3.Range("E15").Value = 4351.950000000001
$ws3.Range("F15").Value = 0.789659255679072

$ws3.Range("D16").Value = 50628.19
$ws3.Range("E16").Value = 4093.040000000001
$ws3.Range("F16").Value = 0.9252019737129447

$ws3.Range("D19").Value = 85738.79000000001
$ws3.Range("E19").Value = 19474.08
$ws3.Range("F19").Value = 0.814907814984992
